# Apply "Updated symbol list" edits to the crypto price sheet.
# The Price column (D) holds numeric-looking values stored as text, so we
# temporarily force a Text number format before writing the new value and
# then clear the formatting again so the cell keeps its original (default)
# style while still being saved back out as a text/string value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# Column D (Price) updates
Set-TextValue "D2"  "243.94"
Set-TextValue "D3"  "25.23"
Set-TextValue "D4"  "5.186"
Set-TextValue "D5"  "0.05732"
Set-TextValue "D6"  "6.498"
Set-TextValue "D7"  "3.111"
Set-TextValue "D8"  "0.8090"
Set-TextValue "D9"  "0.8420"
Set-TextValue "D10" "0.1340"
Set-TextValue "D11" "0.06959"
Set-TextValue "D12" "0.02831"
Set-TextValue "D13" "0.09361"
Set-TextValue "D14" "0.001508"
Set-TextValue "D16" "0.006157"
Set-TextValue "D20" "0.03130"
Set-TextValue "D22" "3.745"
Set-TextValue "D23" "0.04651"
Set-TextValue "D26" "0.004265"
Set-TextValue "D27" "0.00009697"
Set-TextValue "D40" "0.03610"
Set-TextValue "D41" "0.006318"
Set-TextValue "D42" "0.1051"
Set-TextValue "D44" "0.007336"
Set-TextValue "D45" "0.00005278"
Set-TextValue "D48" "0.002285"

# Column E (Volume(1h)) text updates
$ws.Range("E27").Value = "26NitroExNTXBestin24h"
$ws.Range("E41").Value = "40KickTokenKICK"
